$d = $word.ActiveDocument

# First paragraph: title "عقد كراء سيارة" followed by a placeholder run
# for a "{#logo}" tag (plus leading spacer spaces). The placeholder text
# is being dropped and the paragraph re-centered.
$p1 = $d.Paragraphs.Item(1)

# Center the paragraph (was right-aligned).
$p1.Alignment = 1

# Remove everything in the paragraph after the title run (the spacer
# run, the "{" run, "#" run, "logo" run and the "} " run), leaving just
# the title text and the paragraph mark.
$titleText = "عقد كراء سيارة"
$start = $p1.Range.Start
$paraEnd = $p1.Range.End

$tail = $d.Range($start + $titleText.Length, $paraEnd - 1)
if ($tail.Start -lt $tail.End) {
    $tail.Delete()
}
